# Auto-generated script applying scheduled price/profit refresh to Garuda_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 23254.25
$ws.Range("I21").Value = 17672.334
$ws.Range("K21").Value = 17672.334
$ws.Range("M21").Value = -17204.334
$ws.Range("H23").Value = 23254.25
$ws.Range("I23").Value = 17672.334
$ws.Range("K23").Value = 17672.334
$ws.Range("M23").Value = -17438.334
$ws.Range("H29").Value = 2643.75
$ws.Range("J29").Value = 3500
$ws.Range("L29").Value = 10500
$ws.Range("N29").Value = -11062
$ws.Range("H38").Value = 464
$ws.Range("I38").Value = 196
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 588
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = -216
$ws.Range("N38").Value = -3744
$ws.Range("H39").Value = 237.19048
$ws.Range("I39").Value = 34.727272
$ws.Range("J39").Value = 459.9
$ws.Range("K39").Value = 104.181816
$ws.Range("L39").Value = 1379.7
$ws.Range("M39").Value = 191.818184
$ws.Range("N39").Value = -1971.7
$ws.Range("H58").Value = 4056.1853
$ws.Range("I58").Value = 353.2
$ws.Range("J58").Value = 6234.4116
$ws.Range("K58").Value = 1059.6
$ws.Range("L58").Value = 18703.2348
$ws.Range("M58").Value = -909.5999999999999
$ws.Range("N58").Value = -19003.2348
$ws.Range("H62").Value = 3122
$ws.Range("I62").Value = 3398.8
$ws.Range("K62").Value = 3398.8
$ws.Range("M62").Value = -2774.8
$ws.Range("H65").Value = 3122
$ws.Range("I65").Value = 3398.8
$ws.Range("K65").Value = 16994
$ws.Range("M65").Value = -13874
$ws.Range("H87").Value = 32800
$ws.Range("J87").Value = 34666.668
$ws.Range("L87").Value = 34666.668
$ws.Range("N87").Value = -37162.668
$ws.Range("H90").Value = 32800
$ws.Range("J90").Value = 34666.668
$ws.Range("L90").Value = 104000.004
$ws.Range("N90").Value = -116480.004
$ws.Range("H118").Value = 768.3333
$ws.Range("I118").Value = 656.36365
$ws.Range("J118").Value = 2000
$ws.Range("K118").Value = 1969.09095
$ws.Range("L118").Value = 6000
$ws.Range("M118").Value = -312.09095
$ws.Range("N118").Value = -9314
$ws.Range("H125").Value = 1489.3334
$ws.Range("J125").Value = 1759
$ws.Range("L125").Value = 15831
$ws.Range("N125").Value = -20751
$ws.Range("H132").Value = 2269188.2
$ws.Range("I132").Value = 2343481.2
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 7030443.600000001
$ws.Range("L132").Value = 9750
$ws.Range("M132").Value = -7027913.600000001
$ws.Range("N132").Value = -14810
$ws.Range("H137").Value = 22729922
$ws.Range("I137").Value = 1110
$ws.Range("K137").Value = 3330
$ws.Range("M137").Value = -780
$ws.Range("H138").Value = 2943.44
$ws.Range("J138").Value = 3650.845
$ws.Range("L138").Value = 10952.535
$ws.Range("N138").Value = -21232.535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 36991.188
$ws.Range("I32").Value = 35991.08
$ws.Range("K32").Value = 35991.08
$ws.Range("M32").Value = -35704.08
$ws.Range("H132").Value = 5891.518
$ws.Range("I132").Value = 6899.25
$ws.Range("J132").Value = 3372.1875
$ws.Range("K132").Value = 20697.75
$ws.Range("L132").Value = 10116.5625
$ws.Range("M132").Value = -18167.75
$ws.Range("N132").Value = -15176.5625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 16656
$ws.Range("J50").Value = 16656
$ws.Range("L50").Value = 16656
$ws.Range("N50").Value = -17906
$ws.Range("H51").Value = 14867
$ws.Range("J51").Value = 22645
$ws.Range("L51").Value = 22645
$ws.Range("N51").Value = -24117
$ws.Range("H60").Value = 23146.666
$ws.Range("J60").Value = 27376
$ws.Range("L60").Value = 27376
$ws.Range("N60").Value = -28398
$ws.Range("H61").Value = 14867
$ws.Range("J61").Value = 22645
$ws.Range("L61").Value = 22645
$ws.Range("N61").Value = -23341
$ws.Range("H132").Value = 3789679.5
$ws.Range("I132").Value = 1542
$ws.Range("K132").Value = 4626
$ws.Range("M132").Value = -2096

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1164.7142
$ws.Range("J34").Value = 1380.6
$ws.Range("L34").Value = 4141.799999999999
$ws.Range("N34").Value = -4309.799999999999
$ws.Range("H39").Value = 3840
$ws.Range("J39").Value = 3840
$ws.Range("L39").Value = 11520
$ws.Range("N39").Value = -12108
$ws.Range("H55").Value = 2223.3333
$ws.Range("J55").Value = 2223.3333
$ws.Range("L55").Value = 6669.999899999999
$ws.Range("N55").Value = -7023.999899999999
$ws.Range("H131").Value = 736.12
$ws.Range("I131").Value = 298.57144
$ws.Range("J131").Value = 807.3488
$ws.Range("K131").Value = 895.71432
$ws.Range("L131").Value = 2422.0464
$ws.Range("M131").Value = 4144.28568
$ws.Range("N131").Value = -12502.0464
$ws.Range("H133").Value = 2327.7666
$ws.Range("I133").Value = 2707.3684
$ws.Range("J133").Value = 1672.091
$ws.Range("K133").Value = 8122.1052
$ws.Range("L133").Value = 5016.272999999999
$ws.Range("M133").Value = -3062.1052
$ws.Range("N133").Value = -15136.273
$ws.Range("H134").Value = 1659.68
$ws.Range("I134").Value = 1376.7894
$ws.Range("J134").Value = 2555.5
$ws.Range("K134").Value = 4130.3682
$ws.Range("L134").Value = 7666.5
$ws.Range("M134").Value = 939.6318000000001
$ws.Range("N134").Value = -17806.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 15874434
$ws.Range("I61").Value = 1312.8572
$ws.Range("J61").Value = 47620676
$ws.Range("K61").Value = 1312.8572
$ws.Range("L61").Value = 47620676
$ws.Range("M61").Value = -1110.8572
$ws.Range("N61").Value = -47621080
$ws.Range("H113").Value = 15874434
$ws.Range("I113").Value = 1312.8572
$ws.Range("J113").Value = 47620676
$ws.Range("K113").Value = 1312.8572
$ws.Range("L113").Value = 47620676
$ws.Range("M113").Value = 857.1428000000001
$ws.Range("N113").Value = -47625016
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H122").Value = 3900.3
$ws.Range("I122").Value = 4820.6
$ws.Range("J122").Value = 2980
$ws.Range("K122").Value = 14461.8
$ws.Range("L122").Value = 8940
$ws.Range("M122").Value = -12011.8
$ws.Range("N122").Value = -13840
$ws.Range("H132").Value = 5907.367
$ws.Range("I132").Value = 6101.4316
$ws.Range("J132").Value = 4199.6
$ws.Range("K132").Value = 18304.2948
$ws.Range("L132").Value = 12598.8
$ws.Range("M132").Value = -15774.2948
$ws.Range("N132").Value = -17658.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 23961
$ws.Range("J86").Value = 23961
$ws.Range("L86").Value = 23961
$ws.Range("N86").Value = -26207
$ws.Range("H89").Value = 23961
$ws.Range("J89").Value = 23961
$ws.Range("L89").Value = 119805
$ws.Range("N89").Value = -131037
$ws.Range("H107").Value = 24123416
$ws.Range("I107").Value = 10417280
$ws.Range("J107").Value = 47619650
$ws.Range("K107").Value = 31251840
$ws.Range("L107").Value = 142858950
$ws.Range("M107").Value = -31249920
$ws.Range("N107").Value = -142862790
$ws.Range("H110").Value = 16000
$ws.Range("J110").Value = 16000
$ws.Range("L110").Value = 16000
$ws.Range("N110").Value = -24180
$ws.Range("H116").Value = 59680
$ws.Range("J116").Value = 59680
$ws.Range("L116").Value = 59680
$ws.Range("N116").Value = -68858
$ws.Range("H133").Value = 43804.08
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 43804.08
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 43804.08
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -53924.08
